$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Insert three new rows just above the current row 25 (which holds the
# "obo" entry). This pushes the existing rows 25-27 down to 28-30 and
# opens up rows 25-27 for the new "cttv-experiment: biological activity"
# entries (up / down / unknown).
# ---------------------------------------------------------------------
$ws.Range("A25:A27").EntireRow.Insert()

# Copy the formatting (styles) from the row-11 pattern (A s=2, B s=4,
# C s=4, D s=10 -- a blank sub-row of an existing multi-row entry) onto
# the three freshly inserted rows so they match the look of the rest of
# the table instead of picking up a generic default style.
$ws.Range("A11:D11").Copy()
$ws.Range("A25:D27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# New cell values (D25/D26/D27 are entered before C25 so new shared
# strings are appended in the same order as the authored workbook).
# ---------------------------------------------------------------------
$ws.Range("D25").Value = "cttvexp:up"
$ws.Range("D26").Value = "cttvexp:down"
$ws.Range("D27").Value = "cttvexp:unknown"
$ws.Range("C25").Value = "biological_subject{properties}{activity}"

# ---------------------------------------------------------------------
# Expand Table1 (the structured table) so it covers the three new rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A5:D30"))

# ---------------------------------------------------------------------
# The row insert does not keep the Hyperlinks collection's ranges in
# sync, so rebuild all six hyperlinks pointing at their correct
# (possibly shifted) cells.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "http://rdf.ebi.ac.uk/resource/ensembl/")
$ws.Hyperlinks.Add($ws.Range("B8"), "http://purl.uniprot.org/uniprot/")
$ws.Hyperlinks.Add($ws.Range("B13"), "http://www.targetvalidation.org/cttv_core/experiment/")
$ws.Hyperlinks.Add($ws.Range("B10"), "http://www.targetvalidation.org/cttv_core")
$ws.Hyperlinks.Add($ws.Range("B28"), "http://purl.obolibrary.org/obo/")
$ws.Hyperlinks.Add($ws.Range("B29"), "http://www.ebi.ac.uk/efo/")

# ---------------------------------------------------------------------
# Keep the _xlnm._FilterDatabase defined name (driven by the table/
# autofilter) aligned with the new table extent.
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$D`$30"
    }
}

# ---------------------------------------------------------------------
# Match the final selection shown in the edited workbook.
# ---------------------------------------------------------------------
$ws.Range("D25").Select() | Out-Null
